# ===========================================================
# Weekly CompStat update for 40th Precinct (crime data refresh)
# ===========================================================
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number + report week dates) ---
$ws.Range("A8").Value = "Volume 31   Number  33"
$ws.Range("C9").Value = "Report Covering the Week  8/12/2024  Through  8/18/2024"

# --- Plain numeric value updates ---
$ws.Range("N14").Value = -90.74074074074
$ws.Range("E15").Value = -100
$ws.Range("J15").Value = 30
$ws.Range("K15").Value = -3.333333333333
$ws.Range("L15").Value = 52.631578947368
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 19
$ws.Range("E16").Value = -68.421052631578
$ws.Range("F16").Value = 48
$ws.Range("G16").Value = 61
$ws.Range("H16").Value = -21.311475409836
$ws.Range("I16").Value = 402
$ws.Range("J16").Value = 402
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 8.648648648648
$ws.Range("M16").Value = 37.671232876712
$ws.Range("N16").Value = -65.90330788804
$ws.Range("C17").Value = 15
$ws.Range("D17").Value = 19
$ws.Range("E17").Value = -21.052631578947
$ws.Range("F17").Value = 79
$ws.Range("G17").Value = 77
$ws.Range("H17").Value = 2.597402597402
$ws.Range("I17").Value = 648
$ws.Range("J17").Value = 620
$ws.Range("K17").Value = 4.516129032258
$ws.Range("L17").Value = 20.222634508348
$ws.Range("M17").Value = 139.114391143911
$ws.Range("N17").Value = -9.623430962343
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 8
$ws.Range("E18").Value = -25
$ws.Range("F18").Value = 27
$ws.Range("G18").Value = 25
$ws.Range("H18").Value = 8
$ws.Range("I18").Value = 217
$ws.Range("J18").Value = 188
$ws.Range("K18").Value = 15.425531914893
$ws.Range("L18").Value = -2.690582959641
$ws.Range("M18").Value = 65.648854961832
$ws.Range("N18").Value = -75.424688561721
$ws.Range("C19").Value = 35
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = 150
$ws.Range("F19").Value = 90
$ws.Range("G19").Value = 63
$ws.Range("H19").Value = 42.857142857142
$ws.Range("I19").Value = 582
$ws.Range("J19").Value = 452
$ws.Range("K19").Value = 28.761061946902
$ws.Range("L19").Value = 26.247288503253
$ws.Range("M19").Value = 124.710424710425
$ws.Range("N19").Value = 25.701943844492
$ws.Range("C20").Value = 9
$ws.Range("D20").Value = 9
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 31
$ws.Range("G20").Value = 30
$ws.Range("H20").Value = 3.333333333333
$ws.Range("I20").Value = 147
$ws.Range("J20").Value = 233
$ws.Range("K20").Value = -36.909871244635
$ws.Range("L20").Value = -28.640776699029
$ws.Range("M20").Value = 81.481481481481
$ws.Range("N20").Value = -66.590909090909
$ws.Range("C21").Value = 71
$ws.Range("D21").Value = 70
$ws.Range("E21").Value = 1.428571428571
$ws.Range("F21").Value = 279
$ws.Range("G21").Value = 258
$ws.Range("H21").Value = 8.13953488372
$ws.Range("I21").Value = 2030
$ws.Range("J21").Value = 1930
$ws.Range("K21").Value = 5.181347150259
$ws.Range("L21").Value = 11.17196056955
$ws.Range("M21").Value = 91.871455576559
$ws.Range("N21").Value = -46.225165562913
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 100
$ws.Range("F22").Value = 10
$ws.Range("G22").Value = 10
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 54
$ws.Range("J22").Value = 42
$ws.Range("K22").Value = 28.571428571428
$ws.Range("L22").Value = -11.475409836065
$ws.Range("M22").Value = 38.461538461538
$ws.Range("C23").Value = 8
$ws.Range("D23").Value = 14
$ws.Range("E23").Value = -42.857142857142
$ws.Range("F23").Value = 39
$ws.Range("G23").Value = 37
$ws.Range("H23").Value = 5.405405405405
$ws.Range("I23").Value = 298
$ws.Range("J23").Value = 320
$ws.Range("K23").Value = -6.875
$ws.Range("L23").Value = 14.615384615384
$ws.Range("M23").Value = 59.358288770053
$ws.Range("C24").Value = 37
$ws.Range("D24").Value = 40
$ws.Range("E24").Value = -7.5
$ws.Range("F24").Value = 152
$ws.Range("G24").Value = 159
$ws.Range("H24").Value = -4.40251572327
$ws.Range("I24").Value = 952
$ws.Range("J24").Value = 1038
$ws.Range("K24").Value = -8.285163776493
$ws.Range("L24").Value = -11.851851851851
$ws.Range("M24").Value = 14.560770156438
$ws.Range("C25").Value = 19
$ws.Range("E25").Value = 58.333333333333
$ws.Range("F25").Value = 77
$ws.Range("G25").Value = 56
$ws.Range("H25").Value = 37.5
$ws.Range("I25").Value = 397
$ws.Range("J25").Value = 393
$ws.Range("K25").Value = 1.017811704834
$ws.Range("L25").Value = -27.422303473491
$ws.Range("C26").Value = 25
$ws.Range("D26").Value = 29
$ws.Range("E26").Value = -13.793103448275
$ws.Range("F26").Value = 103
$ws.Range("G26").Value = 104
$ws.Range("H26").Value = -0.961538461538
$ws.Range("I26").Value = 795
$ws.Range("J26").Value = 719
$ws.Range("K26").Value = 10.570236439499
$ws.Range("L26").Value = 21.745788667687
$ws.Range("M26").Value = 9.053497942386
$ws.Range("F27").Value = 6
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 40
$ws.Range("J27").Value = 42
$ws.Range("K27").Value = -4.761904761904
$ws.Range("L27").Value = 29.032258064516
$ws.Range("C28").Value = 2
$ws.Range("F28").Value = 12
$ws.Range("G28").Value = 4
$ws.Range("I28").Value = 107
$ws.Range("K28").Value = 44.594594594594
$ws.Range("L28").Value = 98.148148148148
$ws.Range("C29").Value = 1
$ws.Range("I29").Value = 32
$ws.Range("K29").Value = 68.421052631578
$ws.Range("L29").Value = -25.581395348837
$ws.Range("M29").Value = -8.571428571428
$ws.Range("N29").Value = -79.084967320261
$ws.Range("C30").Value = 1
$ws.Range("I30").Value = 27
$ws.Range("K30").Value = 50
$ws.Range("L30").Value = -20.588235294117
$ws.Range("M30").Value = -3.571428571428
$ws.Range("N30").Value = -79.850746268656
$ws.Range("I33").Value = 2
$ws.Range("K33").Value = -50
$ws.Range("L33").Value = -71.428571428571

# --- Cells that change data type (number <-> placeholder text) ---
# Helper: force a literal text value (apostrophe prefix), then fix up the
# cell style by pasting formats from a reference cell that already carries
# the desired style in this workbook, so number formats / quote-prefix noise
# do not leak a brand new style index into the sheet.
function Set-TextCell($addr, $text, $styleRef) {
    $dst = $ws.Range($addr)
    $dst.Value = "" + $text
    $ws.Range($styleRef).Copy() | Out-Null
    $dst.PasteSpecial(-4122) | Out-Null
}
function Set-NumberCell($addr, $value, $styleRef) {
    $dst = $ws.Range($addr)
    $dst.Value = $value
    $ws.Range($styleRef).Copy() | Out-Null
    $dst.PasteSpecial(-4122) | Out-Null
}

Set-TextCell "C15" "0" "D14"
Set-TextCell "D28" "0" "D14"
Set-TextCell "E28" "***.*" "N22"
Set-TextCell "G29" "0" "D14"
Set-TextCell "H29" "***.*" "N22"
Set-TextCell "G30" "0" "D14"
Set-TextCell "H30" "***.*" "N22"
Set-NumberCell "F33" 1 "I33"

$ws.Application.CutCopyMode = $false

# --- Column E width: recalculated bestFit width after E15 widened to -100 ---
$ws.Columns("E:E").ColumnWidth = 6.714285714285714

Write-Host "Weekly CompStat figures refreshed."